# Fruta / hortaliza, semanal
# Insert two new weekly rows for "Feria Lagunitas de Puerto Montt - Uva" right above
# the current row 128 (pushing the existing data down by two rows), and populate
# them with the new week's observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 128; formatting (incl. date number format) is
# inherited from the row above, same as it would be in the interactive UI.
$ws.Rows("128:129").Insert()

# New row 128: Red Globe, $/caja 20 kilos
$ws.Range("A128").Value = 4
$ws.Range("B128").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C128").Value = "Los Lagos"
$ws.Range("D128").Value = 44603
$ws.Range("E128").Value = 10
$ws.Range("F128").Value = "Fruta"
$ws.Range("G128").Value = 100109
$ws.Range("H128").Value = "Uva"
$ws.Range("I128").Value = 100109001
$ws.Range("J128").Value = "Uva"
$ws.Range("K128").Value = "Red Globe"
$ws.Range("L128").Value = "Primera"
$ws.Range("M128").Value = 300
$ws.Range("N128").Value = 18000
$ws.Range("O128").Value = 19000
$ws.Range("P128").Value = 18500
$ws.Range("Q128").Value = "`$/caja 20 kilos"
$ws.Range("R128").Value = "Región de O'Higgins"
$ws.Range("S128").Value = 925
$ws.Range("T128").Value = 20

# New row 129: Superior Seedless, $/caja 20 kilos
$ws.Range("A129").Value = 4
$ws.Range("B129").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C129").Value = "Los Lagos"
$ws.Range("D129").Value = 44603
$ws.Range("E129").Value = 10
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100109
$ws.Range("H129").Value = "Uva"
$ws.Range("I129").Value = 100109001
$ws.Range("J129").Value = "Uva"
$ws.Range("K129").Value = "Superior Seedless"
$ws.Range("L129").Value = "Primera"
$ws.Range("M129").Value = 300
$ws.Range("N129").Value = 17000
$ws.Range("O129").Value = 18000
$ws.Range("P129").Value = 17500
$ws.Range("Q129").Value = "`$/caja 20 kilos"
$ws.Range("R129").Value = "Región de O'Higgins"
$ws.Range("S129").Value = 875
$ws.Range("T129").Value = 20
